# Regenerate save_data: column G ("K") is recomputed to use the new K
# (strikeouts) metric instead of the old "Strike#" value. Write the
# freshly calculated values into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (recalculated s_vals for the K column)
$newK = @{
    2  = 2
    3  = 2
    4  = 8
    5  = 3
    6  = 9
    7  = 1
    8  = 2
    9  = 4
    10 = 1
    11 = 2
    12 = 4
    13 = 3
    14 = 11
    15 = 5
    16 = 2
    17 = 1
    18 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
